# Insert two new columns before column D (shifts D:K -> F:M)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D:E").Insert()

# Copy number/date formatting from the (now-shifted) first historical column
# into the two new columns, per data block, so the new cells pick up the
# same style (date format row vs. numeric rows) as their row.
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the two new columns with the newest two quarters of data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 466700
$ws.Range("E8").Value = 508000
$ws.Range("D9").Value = 388100
$ws.Range("E9").Value = 423900
$ws.Range("D10").Value = 78600
$ws.Range("E10").Value = 84100
$ws.Range("D12").Value = 13000
$ws.Range("E12").Value = 14000
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 200
$ws.Range("E14").Value = 1700
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 427300
$ws.Range("E17").Value = 480300
$ws.Range("D18").Value = 39400
$ws.Range("E18").Value = 27700
$ws.Range("D20").Value = -4900
$ws.Range("E20").Value = -2500
$ws.Range("D21").Value = 54900
$ws.Range("E21").Value = 45500
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 34500
$ws.Range("E23").Value = 25200
$ws.Range("D24").Value = 1500
$ws.Range("E24").Value = 6800
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 33000
$ws.Range("E26").Value = 18500
$ws.Range("D27").Value = 33000
$ws.Range("E27").Value = 18500
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = -6100
$ws.Range("E29").Value = 3700
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 4900
$ws.Range("E32").Value = 2500
$ws.Range("D33").Value = 26900
$ws.Range("E33").Value = 22200
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 26900
$ws.Range("E35").Value = 22200
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 300200
$ws.Range("E41").Value = 274000
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 280000
$ws.Range("E43").Value = 323100
$ws.Range("D44").Value = 200200
$ws.Range("E44").Value = 195200
$ws.Range("D45").Value = 22100
$ws.Range("E45").Value = 22800
$ws.Range("D46").Value = 802500
$ws.Range("E46").Value = 815000
$ws.Range("D47").Value = 25100
$ws.Range("E47").Value = 28100
$ws.Range("D48").Value = 608900
$ws.Range("E48").Value = 598700
$ws.Range("D49").Value = 37200
$ws.Range("E49").Value = 37500
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 11000
$ws.Range("E52").Value = 13000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1484700
$ws.Range("E54").Value = 1492200
$ws.Range("D57").Value = 206000
$ws.Range("E57").Value = 209300
$ws.Range("D58").Value = 37100
$ws.Range("E58").Value = 23700
$ws.Range("D59").Value = 95600
$ws.Range("E59").Value = 93900
$ws.Range("D60").Value = 338600
$ws.Range("E60").Value = 326900
$ws.Range("D61").Value = 239000
$ws.Range("E61").Value = 262600
$ws.Range("D62").Value = 122500
$ws.Range("E62").Value = 124400
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 700900
$ws.Range("E66").Value = 714600
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 813400
$ws.Range("E72").Value = 792100
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 783800
$ws.Range("E76").Value = 777700
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 26900
$ws.Range("E81").Value = 22200
$ws.Range("D83").Value = 20300
$ws.Range("E83").Value = 20300
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 75000
$ws.Range("E89").Value = 42400
$ws.Range("D91").Value = -23800
$ws.Range("E91").Value = -19200
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -25100
$ws.Range("E94").Value = -19200
$ws.Range("D96").Value = -5600
$ws.Range("E96").Value = -5100
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -17600
$ws.Range("E100").Value = -6900
$ws.Range("D101").Value = -6000
$ws.Range("E101").Value = 1000
$ws.Range("D102").Value = 26200
$ws.Range("E102").Value = 17300

# A handful of historical figures were restated; fix them up in their new (shifted) positions
$ws.Range("H9").Value = 393000
$ws.Range("H10").Value = 80800
$ws.Range("H12").Value = 13400
$ws.Range("H17").Value = 442400
$ws.Range("I17").Value = 457400
$ws.Range("H18").Value = 31400
$ws.Range("I18").Value = 30400
$ws.Range("H20").Value = -2800
$ws.Range("I20").Value = -1100
$ws.Range("H24").Value = 18600
$ws.Range("H26").Value = 9900
$ws.Range("H27").Value = 9900
$ws.Range("H32").Value = 2800
$ws.Range("I32").Value = 1100
$ws.Range("H33").Value = -5000
$ws.Range("H35").Value = -5000
$ws.Range("H81").Value = -5000
